# Add the white paper of AI
#
# - Adds a new row (A4) on the first sheet ("博士 JF westbrook") containing a
#   new shared string about the white paper / entrepreneurship application info.
# - Makes the first sheet the active/selected tab (with cell D8 selected),
#   which also clears the "tabSelected" flag from the sheet that was
#   previously active ("星星", selection stays at G9).

$wb = $excel.ActiveWorkbook

$wsFirst = $wb.Worksheets.Item(1)

# Write the new shared-string cell into row 4, column A.
$wsFirst.Range("A4").Value = "同济的创业的申请每年4月份、10月份各有一次。如果答辩通过了提供资金和场所"

# Make this sheet the active tab, and select D8 on it (matches the diff's
# new <sheetView tabSelected="1" ...><selection activeCell="D8" .../>).
$null = $wsFirst.Activate()
$null = $wsFirst.Range("D8").Select()
